# Update "想去人数" (column F) counts on sheets 展览, 演出, 全部类型
# to match the values output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 327
$ws.Range("F3").Value = 279
$ws.Range("F4").Value = 1206
$ws.Range("F9").Value = 136
$ws.Range("F10").Value = 3412
$ws.Range("F11").Value = 119
$ws.Range("F12").Value = 81
$ws.Range("F13").Value = 64
$ws.Range("F14").Value = 37
$ws.Range("F15").Value = 52
$ws.Range("F16").Value = 581
$ws.Range("F17").Value = 72
$ws.Range("F18").Value = 695
$ws.Range("F19").Value = 202
$ws.Range("F20").Value = 113
$ws.Range("F21").Value = 56
$ws.Range("F22").Value = 52
$ws.Range("F24").Value = 2491
$ws.Range("F25").Value = 5003
$ws.Range("F28").Value = 474
$ws.Range("F29").Value = 1284
$ws.Range("F31").Value = 2206
$ws.Range("F32").Value = 571
$ws.Range("F34").Value = 78
$ws.Range("F35").Value = 94
$ws.Range("F37").Value = 307
$ws.Range("F38").Value = 453
$ws.Range("F41").Value = 450

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 68

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 327
$ws.Range("F3").Value = 279
$ws.Range("F4").Value = 1206
$ws.Range("F9").Value = 136
$ws.Range("F10").Value = 3412
$ws.Range("F11").Value = 119
$ws.Range("F12").Value = 81
$ws.Range("F13").Value = 64
$ws.Range("F14").Value = 68
$ws.Range("F15").Value = 37
$ws.Range("F16").Value = 52
$ws.Range("F17").Value = 581
$ws.Range("F18").Value = 72
$ws.Range("F19").Value = 695
$ws.Range("F20").Value = 202
$ws.Range("F21").Value = 113
$ws.Range("F22").Value = 56
$ws.Range("F23").Value = 52
$ws.Range("F25").Value = 2491
$ws.Range("F26").Value = 5003
$ws.Range("F29").Value = 474
$ws.Range("F30").Value = 1284
$ws.Range("F32").Value = 2206
$ws.Range("F33").Value = 571
$ws.Range("F35").Value = 78
$ws.Range("F36").Value = 94
$ws.Range("F38").Value = 307
$ws.Range("F39").Value = 453
$ws.Range("F42").Value = 450
